# Re-orders the "Recorded By" (column G) comma-separated list of names on
# the "Session Analysis Results" sheet so that a literal "System" entry is
# moved to the end of the list, instead of wherever it happened to sort
# before (typically first). Other entries keep their relative order.
#
# Example: "System, dnasr281@gmail.com"        -> "dnasr281@gmail.com, System"
#          "system, System, backup@backdoor.com" -> "system, backup@backdoor.com, System"
#
# Rows where "System" is already last, or isn't present at all, are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value()

    if ($null -eq $value) {
        continue
    }

    $text = [string]$value
    $parts = $text -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    # NOTE: -eq/-ne/-contains on strings are case-INsensitive, which would
    # wrongly also strip a lowercase "system" entry. Use the case-sensitive
    # .Equals() instance method instead so only an exact "System" token
    # (capital S) is relocated.
    $hasExactSystem = $false
    $rest = @()
    foreach ($p in $parts) {
        if ($p.Equals('System')) {
            $hasExactSystem = $true
        } else {
            $rest += $p
        }
    }

    if ($hasExactSystem) {
        $newParts = @($rest) + @('System')
        $newText = [string]::Join(', ', $newParts)

        if (-not $newText.Equals($text)) {
            $cell.Value = $newText
        }
    }
}
